# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.691.70"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "2.519.32"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.45"
$ws.Range("E5").Value = "  +3.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.21"
$ws.Range("E6").Value = "  -3.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -1.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.16"
$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.70"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("E13").Value = "  -2.43%  "

$ws.Range("D14").Value = "2.907.48"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("E15").Value = "  +3.99%  "

$ws.Range("D16").Value = "2.514.12"
$ws.Range("E16").Value = "  -4.13%  "

$ws.Range("E17").Value = "  -2.35%  "

$ws.Range("D18").Value = "42.766.49"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.15"
$ws.Range("E19").Value = "  -4.66%  "

$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  -2.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.53"
$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.29"
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.24"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("E24").Value = "  +0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("E25").Value = "  -2.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.17"
$ws.Range("E26").Value = "  -3.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.44%  "

$ws.Range("E28").Value = "  +11.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.26"
$ws.Range("E29").Value = "  +2.92%  "

$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("E31").Value = "  -4.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.93"
$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.85"
$ws.Range("E33").Value = "  +1.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.33"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.09"
$ws.Range("E35").Value = "  -3.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0786"
$ws.Range("E36").Value = "  -2.42%  "

$ws.Range("E37").Value = "  -5.38%  "

$ws.Range("E38").Value = "  -2.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.28"
$ws.Range("E39").Value = "  -6.03%  "

$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("E41").Value = "  +1.04%  "

$ws.Range("E42").Value = "  -1.23%  "

$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("D44").Value = "2.065.99"
$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0302"
$ws.Range("E46").Value = "  -1.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.01"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.81"
$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.66"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").Value = "2.762.73"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.190"
$ws.Range("E51").Value = "  -1.04%  "
